# Update the "EXECUTABLE" column (column E) on the "Master" sheet:
# rows 2-14 change from "YES" to "NO" (row 15 and below remain "YES").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Master")

$ws.Range("E2:E14").Value = "NO"
